$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weight/group table for the ceramide scoring sheet.
# Column A = Type (shared string label), Column B = Weight, Column C = Group
$rows = @(
    @("FA1_[LCB+H]+",               20.5, 1),
    @("FA1_[LCB-H2O+H]+",           20.5, 1),
    @("FA1_[LCB-2xH2O+H]+",         20.5, 1),
    @("FA1_[LCB-3xH2O+H]+",          2,   1),
    @("FA1_[LCB-H2O-CH2O+H]+",       0.5, 1),
    @("FA1_[LCB-2xH2O-CH2O+H]+",     0.5, 1),
    @("FA2_[FA-HO+NH3]+",            8,   3),
    @("[M-H2O+H]+",                 20.5, 2),
    @("[M-2xH2O+H]+",                4,   2),
    @("[M-3xH2O+H]+",                2,   2),
    @("[M-H2O-CH2O+H]+",             0.5, 2),
    @("[M-2xH2O-CH2O+H]+",           0.5, 2)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("B13").Select()
